# Update for new trial protocol and pis.
#
# 1. Refresh the "today" date placeholder text (datetimeFigureOut field)
#    across the slide master, every slide layout, and the notes master.
# 2. Slide 5 ("UC Flare" schematic): shift the "UC Flare" textbox and the
#    dotted green connector under it to the right, shrink the connector,
#    and shrink / retitle the "Screening (Max 2 weeks)" textbox to
#    "Screening (4 weeks)".

$p = $ppt.ActivePresentation
$EMU_PER_PT = 12700.0
# tiny nudge so that EMU -> pt -> EMU round trips land back on the exact
# integer EMU value instead of being truncated a hair short.
$EPS = 0.5 / $EMU_PER_PT

function ToPt($emu) {
    return ($emu / $EMU_PER_PT) + $EPS
}

function Set-DatePlaceholderText($shapes, $newText) {
    foreach ($sh in $shapes) {
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "8/12/25"

# --- Slide master ---
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# --- Every slide layout ---
foreach ($layout in $p.SlideMaster.CustomLayouts) {
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# --- Notes master ---
Set-DatePlaceholderText $p.NotesMaster.Shapes $newDate

# --- Slide 5 shape tweaks ---
$slide5 = $p.Slides.Item(5)

# "UC Flare" textbox: slide right, same size.
$ucFlare = $null
foreach ($sh in $slide5.Shapes) {
    if ($sh.Name -eq "TextBox 9") { $ucFlare = $sh }
}
$ucFlare.Left = ToPt 1605670

# Dotted green connector under "UC Flare": move right, shorten, and drop
# the slight vertical flip/slope it used to have (now perfectly flat).
$connector16 = $null
foreach ($sh in $slide5.Shapes) {
    if ($sh.Name -eq "Straight Connector 16") { $connector16 = $sh }
}
$connector16.LockAspectRatio = $true
$connector16.VerticalFlip = $false
$connector16.Left = ToPt 1692876
$connector16.Top = ToPt 3080064
$connector16.Width = ToPt 1314811
$connector16.Height = ToPt 0

# "Screening / (Max 2 weeks)" textbox: slide right, narrower, and the
# second line of text now reads "(4 weeks)".
$screening = $null
foreach ($sh in $slide5.Shapes) {
    if ($sh.Name -eq "TextBox 25") { $screening = $sh }
}
$screening.Left = ToPt 1605670
$screening.Width = ToPt 1043876

$tr = $screening.TextFrame.TextRange
$fullText = $tr.Text
$idx = $fullText.IndexOf("(Max 2 weeks)")
$target = $tr.Characters($idx + 1, 13)
$target.Text = "(4 weeks)"
